$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 119114.58
$ws.Range("I40").Value = 1500590
$ws.Range("J40").Value = 3991.625
$ws.Range("K40").Value = 1500590
$ws.Range("L40").Value = 3991.625
$ws.Range("M40").Value = -1500415
$ws.Range("N40").Value = -4341.625
$ws.Range("H74").Value = 9398
$ws.Range("I74").Value = 8996.666999999999
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 8996.666999999999
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -8060.666999999999
$ws.Range("N74").Value = -11872
$ws.Range("H77").Value = 9398
$ws.Range("I77").Value = 8996.666999999999
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 44983.335
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -40303.335
$ws.Range("N77").Value = -59360
$ws.Range("H125").Value = 3830.0625
$ws.Range("I125").Value = 1528
$ws.Range("J125").Value = 5620.5557
$ws.Range("K125").Value = 13752
$ws.Range("L125").Value = 50585.0013
$ws.Range("M125").Value = -11292
$ws.Range("N125").Value = -55505.0013
$ws.Range("H135").Value = 394.68
$ws.Range("I135").Value = 414.69565
$ws.Range("J135").Value = 164.5
$ws.Range("K135").Value = 3732.26085
$ws.Range("L135").Value = 1480.5
$ws.Range("M135").Value = -1197.26085
$ws.Range("N135").Value = -6550.5
$ws.Range("H138").Value = 3687.16
$ws.Range("J138").Value = 3762.1458
$ws.Range("L138").Value = 11286.4374
$ws.Range("N138").Value = -21566.4374

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2927.7878
$ws.Range("I122").Value = 2184.5652
$ws.Range("K122").Value = 6553.6956
$ws.Range("M122").Value = -4103.6956
$ws.Range("H131").Value = 99592.60000000001
$ws.Range("J131").Value = 99592.60000000001
$ws.Range("L131").Value = 99592.60000000001
$ws.Range("N131").Value = -109672.6
$ws.Range("H132").Value = 2388.853
$ws.Range("I132").Value = 2107.3667
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 6322.1001
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -3792.1001
$ws.Range("N132").Value = -18560

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2138.5264
$ws.Range("I20").Value = 2094.2144
$ws.Range("J20").Value = 2262.6
$ws.Range("K20").Value = 2094.2144
$ws.Range("L20").Value = 2262.6
$ws.Range("M20").Value = -1847.2144
$ws.Range("N20").Value = -2756.6
$ws.Range("H86").Value = 2224.818
$ws.Range("I86").Value = 1140.25
$ws.Range("J86").Value = 5117
$ws.Range("K86").Value = 1140.25
$ws.Range("L86").Value = 5117
$ws.Range("M86").Value = -17.25
$ws.Range("N86").Value = -7363
$ws.Range("H89").Value = 2224.818
$ws.Range("I89").Value = 1140.25
$ws.Range("J89").Value = 5117
$ws.Range("K89").Value = 5701.25
$ws.Range("L89").Value = 25585
$ws.Range("M89").Value = -85.25
$ws.Range("N89").Value = -36817
$ws.Range("H107").Value = 3060.6365
$ws.Range("I107").Value = 2964.111
$ws.Range("K107").Value = 2964.111
$ws.Range("M107").Value = -1044.111
$ws.Range("H130").Value = 96810.336
$ws.Range("J130").Value = 96810.336
$ws.Range("L130").Value = 96810.336
$ws.Range("N130").Value = -106850.336
$ws.Range("H134").Value = 3177582.5
$ws.Range("I134").Value = 3923778.8
$ws.Range("K134").Value = 11771336.4
$ws.Range("M134").Value = -11768801.4

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1625.5
$ws.Range("J31").Value = 1579.86
$ws.Range("L31").Value = 1579.86
$ws.Range("N31").Value = -2169.86
$ws.Range("H34").Value = 1625.5
$ws.Range("J34").Value = 1579.86
$ws.Range("L34").Value = 1579.86
$ws.Range("N34").Value = -1983.86
$ws.Range("H132").Value = 3632.2222
$ws.Range("I132").Value = 3027.4285
$ws.Range("K132").Value = 9082.2855
$ws.Range("M132").Value = -6552.2855

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1300
$ws.Range("I3").Value = 1300
$ws.Range("K3").Value = 3900
$ws.Range("M3").Value = -3788
$ws.Range("H39").Value = 4576
$ws.Range("J39").Value = 4576
$ws.Range("L39").Value = 13728
$ws.Range("N39").Value = -14316
$ws.Range("H106").Value = 9999.5
$ws.Range("J106").Value = 9999.5
$ws.Range("L106").Value = 29998.5
$ws.Range("N106").Value = -31890.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 23330
$ws.Range("J52").Value = 23330
$ws.Range("L52").Value = 23330
$ws.Range("N52").Value = -23848
$ws.Range("H53").Value = 25000
$ws.Range("I53").Value = 25000
$ws.Range("K53").Value = 25000
$ws.Range("M53").Value = -24369
$ws.Range("H97").Value = 1072.2333
$ws.Range("I97").Value = 913.38464
$ws.Range("K97").Value = 913.38464
$ws.Range("M97").Value = -417.38464
$ws.Range("H123").Value = 29888
$ws.Range("J123").Value = 29888
$ws.Range("L123").Value = 29888
$ws.Range("N123").Value = -34788
$ws.Range("H132").Value = 3611.923
$ws.Range("I132").Value = 3158.5881
$ws.Range("J132").Value = 4468.222
$ws.Range("K132").Value = 9475.764299999999
$ws.Range("L132").Value = 13404.666
$ws.Range("M132").Value = -6945.764299999999
$ws.Range("N132").Value = -18464.666
$ws.Range("H140").Value = 111784.914
$ws.Range("J140").Value = 700001
$ws.Range("L140").Value = 700001
$ws.Range("N140").Value = -710361

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 92814.7
$ws.Range("I93").Value = 130306.86
$ws.Range("J93").Value = 5333
$ws.Range("K93").Value = 130306.86
$ws.Range("L93").Value = 5333
$ws.Range("M93").Value = -129058.86
$ws.Range("N93").Value = -7829
$ws.Range("H132").Value = 10848.786
$ws.Range("I132").Value = 10914.154
$ws.Range("K132").Value = 32742.462
$ws.Range("M132").Value = -30212.462
$ws.Range("H136").Value = 3703.9048
$ws.Range("J136").Value = 3748.5
$ws.Range("L136").Value = 11245.5
$ws.Range("N136").Value = -16345.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2472.625
$ws.Range("I107").Value = 756.4
$ws.Range("K107").Value = 2269.2
$ws.Range("M107").Value = -349.1999999999998
$ws.Range("H126").Value = 3484.9062
$ws.Range("I126").Value = 3300.5925
$ws.Range("K126").Value = 9901.7775
$ws.Range("M126").Value = -7431.7775
$ws.Range("H130").Value = 49995
$ws.Range("J130").Value = 49995
$ws.Range("L130").Value = 49995
$ws.Range("N130").Value = -60035
$ws.Range("H136").Value = 21805.314
$ws.Range("I136").Value = 1846
$ws.Range("J136").Value = 55421
$ws.Range("K136").Value = 5538
$ws.Range("L136").Value = 166263
$ws.Range("M136").Value = -2988
$ws.Range("N136").Value = -171363
